$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update C2 value (Fecha_Proyectada for Proyecto A) 45992 -> 46011
$ws.Range("C2").Value = 46011

# Clear B3 and C3 (dates removed for Proyecto B)
$ws.Range("B3").Value = $null
$ws.Range("C3").Value = $null

# Add new row: Proyecto C
$ws.Range("A2:C2").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Proyecto C"
$ws.Range("B4").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("B4").NumberFormat = "mm-dd-yy"
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Tabla1")
$table.Resize($ws.Range("A1:C4"))

$ws.Range("D19").Select()
